# feat: add 2022-Q3 data
#
# Before: "总计" (summary) + "2021-Q3" (fund holdings for 600165 as of 2021-Q3)
# After:  "总计" (summary, +1 row) + "2022-Q3" (new fund holdings) + "2021-Q3"
#         (unchanged fund holdings, now the 3rd tab)
#
# The new "2022-Q3" sheet is inserted between "总计" and "2021-Q3", and a new
# summary row for it is added to "总计" right above the existing 2021-Q3 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: locate the two existing sheets.
# ---------------------------------------------------------------------------
$wsTotal   = $wb.Worksheets.Item("总计")
$wsOld2021 = $wb.Worksheets.Item("2021-Q3")

# ---------------------------------------------------------------------------
# Step 2: duplicate the current "2021-Q3" sheet, placing the copy right after
# it. The copy will keep holding the original 2021-Q3 fund-holdings data
# (and formatting) untouched -- it becomes the new "2021-Q3" tab.
# The original sheet object is then repurposed/renamed to "2022-Q3" and its
# contents are replaced below, which keeps its original sheetId/r:id (matches
# an in-place rename) while the duplicate gets a fresh sheetId (matches a
# newly-created sheet).
# ---------------------------------------------------------------------------
$wsOld2021.Copy($null, $wsOld2021)
$wsNew2021 = $wb.Worksheets.Item(3)

$ws2022 = $wsOld2021
$ws2022.Name = "2022-Q3"
$wsNew2021.Name = "2021-Q3"

# ---------------------------------------------------------------------------
# Step 3: wipe the (renamed) "2022-Q3" sheet's old data and restyle its
# header row / index column to match "总计"'s header styling.
# ---------------------------------------------------------------------------
$ws2022.Cells.Clear()

$wsTotal.Range("B1").Copy()
$ws2022.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$ws2022.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 4: write the new 2022-Q3 fund-holdings data.
# Numeric-looking values that must keep their exact text (leading/trailing
# zeros, e.g. a fund code or a padded percentage) are entered with a leading
# apostrophe so Excel stores them as text instead of re-normalizing them as
# numbers.
# ---------------------------------------------------------------------------
$ws2022.Range("B1").Value = "基金代码"
$ws2022.Range("C1").Value = "基金名称"
$ws2022.Range("D1").Value = "基金规模"
$ws2022.Range("E1").Value = "股票总仓位"
$ws2022.Range("F1").Value = "仓位占比"
$ws2022.Range("G1").Value = "持有市值(亿元)"
$ws2022.Range("H1").Value = "仓位排名"

$ws2022.Range("A2").Value = 0
$ws2022.Range("B2").Value = "'003456"
$ws2022.Range("C2").Value = "信澳新目标灵活配置混合"
$ws2022.Range("D2").Value = "'0.39"
$ws2022.Range("E2").Value = "'94.17"
$ws2022.Range("F2").Value = "'1.53"
$ws2022.Range("G2").Value = "'0.0060"
$ws2022.Range("H2").Value = 6

# ---------------------------------------------------------------------------
# Step 5: insert a new summary row for 2022-Q3 into "总计", directly above
# the existing 2021-Q3 summary row, and renumber the index column.
# ---------------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.01
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Range("A3").Value = 1

# ---------------------------------------------------------------------------
# Step 6: restore the original active-tab selection (the "2021-Q3" sheet was
# the active/selected tab before the edit).
# ---------------------------------------------------------------------------
$wsNew2021.Activate()
